$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '43.553.57'
Set-TextValue $ws.Range("E2") '  +0.51%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.341.01'
Set-TextValue $ws.Range("E3") '  -1.42%  '

# Row 4
Set-TextValue $ws.Range("E4") '  -0.04%  '

# Row 5
Set-TextValue $ws.Range("D5") '304.98'
Set-TextValue $ws.Range("E5") '  -1.54%  '

# Row 6
Set-TextValue $ws.Range("D6") '101.97'
Set-TextValue $ws.Range("E6") '  -2.55%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.513'
Set-TextValue $ws.Range("E7") '  -2.40%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.514'
Set-TextValue $ws.Range("E9") '  -1.01%  '

# Row 10
Set-TextValue $ws.Range("D10") '35.30'
Set-TextValue $ws.Range("E10") '  -2.78%  '

# Row 11
Set-TextValue $ws.Range("B11") 'Dogecoin'
Set-TextValue $ws.Range("C11") 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range("D11") '0.0799'
Set-TextValue $ws.Range("E11") '  -1.73%  '

# Row 12
Set-TextValue $ws.Range("B12") 'TRON'
Set-TextValue $ws.Range("C12") 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D12") '0.113'
Set-TextValue $ws.Range("E12") '  +0.44%  '

# Row 13
Set-TextValue $ws.Range("B13") 'Polkadot'
Set-TextValue $ws.Range("C13") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D13") '6.82'
Set-TextValue $ws.Range("E13") '  -2.68%  '

# Row 14
Set-TextValue $ws.Range("B14") 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range("C14") 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D14") '2.703.56'
Set-TextValue $ws.Range("E14") '  -1.47%  '

# Row 15
Set-TextValue $ws.Range("B15") 'Chainlink'
Set-TextValue $ws.Range("C15") 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D15") '15.62'
Set-TextValue $ws.Range("E15") '  -0.31%  '

# Row 16
Set-TextValue $ws.Range("B16") 'WrappedEther'
Set-TextValue $ws.Range("C16") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D16") '2.328.15'
Set-TextValue $ws.Range("E16") '  -1.96%  '

# Row 17
Set-TextValue $ws.Range("B17") 'Polygon'
Set-TextValue $ws.Range("C17") 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D17") '0.809'
Set-TextValue $ws.Range("E17") '  -1.28%  '

# Row 18
Set-TextValue $ws.Range("B18") 'WrappedBTC'
Set-TextValue $ws.Range("C18") 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range("D18") '43.465.27'
Set-TextValue $ws.Range("E18") '  +0.34%  '

# Row 19
Set-TextValue $ws.Range("B19") 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range("C19") 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D19") '11.86'
Set-TextValue $ws.Range("E19") '  -1.19%  '

# Row 20
Set-TextValue $ws.Range("B20") 'ShibaInu'
Set-TextValue $ws.Range("C20") 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D20") '0.0₃0909'
Set-TextValue $ws.Range("E20") '  -1.94%  '

# Row 21
Set-TextValue $ws.Range("B21") 'Uniswap'
Set-TextValue $ws.Range("C21") 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D21") '6.13'
Set-TextValue $ws.Range("E21") '  -2.68%  '

# Row 22
Set-TextValue $ws.Range("B22") 'Litecoin'
Set-TextValue $ws.Range("C22") 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D22") '68.37'
Set-TextValue $ws.Range("E22") '  -0.16%  '

# Row 23
Set-TextValue $ws.Range("B23") 'BitcoinCash'
Set-TextValue $ws.Range("C23") 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D23") '238.37'
Set-TextValue $ws.Range("E23") '  -1.55%  '

# Row 24
Set-TextValue $ws.Range("B24") 'ImmutableX'
Set-TextValue $ws.Range("C24") 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D24") '1.99'
Set-TextValue $ws.Range("E24") '  -3.48%  '

# Row 25
Set-TextValue $ws.Range("B25") 'PancakeSwap'
Set-TextValue $ws.Range("C25") 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D25") '2.54'
Set-TextValue $ws.Range("E25") '  -3.05%  '

# Row 26
Set-TextValue $ws.Range("B26") 'Dai'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D26") '1.00'
Set-TextValue $ws.Range("E26") '  -0.31%  '

# Row 27
Set-TextValue $ws.Range("B27") 'EthereumClassic'
Set-TextValue $ws.Range("C27") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D27") '25.07'
Set-TextValue $ws.Range("E27") '  -3.98%  '

# Row 28
Set-TextValue $ws.Range("B28") 'InjectiveProtocol'
Set-TextValue $ws.Range("C28") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D28") '34.68'
Set-TextValue $ws.Range("E28") '  -6.04%  '

# Row 29
Set-TextValue $ws.Range("B29") 'Toncoin'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D29") '2.08'
Set-TextValue $ws.Range("E29") '  -9.45%  '

# Row 30
Set-TextValue $ws.Range("D30") '166.29'
Set-TextValue $ws.Range("E30") '  +2.73%  '

# Row 31
Set-TextValue $ws.Range("B31") 'Cosmos'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D31") '9.27'
Set-TextValue $ws.Range("E31") '  -3.70%  '

# Row 32
Set-TextValue $ws.Range("B32") 'FirstDigitalUSD'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D32") '0.999'
Set-TextValue $ws.Range("E32") '  -0.07%  '

# Row 33
Set-TextValue $ws.Range("B33") 'Filecoin'
Set-TextValue $ws.Range("C33") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D33") '5.08'
Set-TextValue $ws.Range("E33") '  -4.04%  '

# Row 34
Set-TextValue $ws.Range("B34") 'WEMIXToken'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D34") '2.42'
Set-TextValue $ws.Range("E34") '  -4.86%  '

# Row 35
Set-TextValue $ws.Range("B35") 'RenderToken'
Set-TextValue $ws.Range("C35") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D35") '4.52'
Set-TextValue $ws.Range("E35") '  -3.43%  '

# Row 36
Set-TextValue $ws.Range("B36") 'Celestia'
Set-TextValue $ws.Range("C36") 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range("D36") '16.91'
Set-TextValue $ws.Range("E36") '  -7.99%  '

# Row 37
Set-TextValue $ws.Range("B37") 'Hedera'
Set-TextValue $ws.Range("C37") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D37") '0.0707'
Set-TextValue $ws.Range("E37") '  -4.49%  '

# Row 38
Set-TextValue $ws.Range("B38") 'LidoDAOToken'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D38") '2.92'
Set-TextValue $ws.Range("E38") '  -6.81%  '

# Row 39
Set-TextValue $ws.Range("B39") 'ARBITRUM'
Set-TextValue $ws.Range("C39") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D39") '1.83'
Set-TextValue $ws.Range("E39") '  -6.61%  '

# Row 40
Set-TextValue $ws.Range("B40") 'Kaspa'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D40") '0.103'
Set-TextValue $ws.Range("E40") '  -2.87%  '

# Row 41
Set-TextValue $ws.Range("B41") 'Stellar'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D41") '0.112'
Set-TextValue $ws.Range("E41") '  -2.76%  '

# Row 42
Set-TextValue $ws.Range("B42") 'ApeXProtocol'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D42") '2.40'
Set-TextValue $ws.Range("E42") '  -1.44%  '

# Row 43
Set-TextValue $ws.Range("B43") 'Maker'
Set-TextValue $ws.Range("C43") 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D43") '1.994.94'
Set-TextValue $ws.Range("E43") '  -0.51%  '

# Row 44
Set-TextValue $ws.Range("B44") 'VeChain'
Set-TextValue $ws.Range("C44") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D44") '0.0286'
Set-TextValue $ws.Range("E44") '  -2.00%  '

# Row 45
Set-TextValue $ws.Range("B45") 'EnergySwap'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D45") '18.51'
Set-TextValue $ws.Range("E45") '  -9.22%  '

# Row 46
Set-TextValue $ws.Range("B46") 'NEARProtocol'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D46") '2.96'
Set-TextValue $ws.Range("E46") '  -8.01%  '

# Row 47
Set-TextValue $ws.Range("D47") '10.01'
Set-TextValue $ws.Range("E47") '  -3.73%  '

# Row 48
Set-TextValue $ws.Range("B48") 'MultiversX'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue $ws.Range("D48") '56.95'
Set-TextValue $ws.Range("E48") '  -1.84%  '

# Row 49
Set-TextValue $ws.Range("B49") 'THORChain'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D49") '4.89'
Set-TextValue $ws.Range("E49") '  +3.88%  '

# Row 50
Set-TextValue $ws.Range("B50") 'RocketPoolETH'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range("D50") '2.566.42'
Set-TextValue $ws.Range("E50") '  +0.25%  '

# Row 51
Set-TextValue $ws.Range("B51") 'Stacks'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D51") '1.55'
Set-TextValue $ws.Range("E51") '  -1.11%  '
